$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 66
$ws1.Range("F3").Value = 655
$ws1.Range("F4").Value = 228
$ws1.Range("F6").Value = 9901
$ws1.Range("F7").Value = 895
$ws1.Range("F9").Value = 1240
$ws1.Range("F10").Value = 3958
$ws1.Range("F11").Value = 1
$ws1.Range("F12").Value = 5
$ws1.Range("F15").Value = 53
$ws1.Range("F18").Value = 563
$ws1.Range("F21").Value = 1475

# Sheet "演出" (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 66
$ws4.Range("F3").Value = 19
$ws4.Range("F4").Value = 655
$ws4.Range("F5").Value = 228
$ws4.Range("F7").Value = 9901
$ws4.Range("F8").Value = 895
$ws4.Range("F10").Value = 1240
$ws4.Range("F11").Value = 3958
$ws4.Range("F12").Value = 1
$ws4.Range("F13").Value = 5
$ws4.Range("F16").Value = 53
$ws4.Range("F19").Value = 563
$ws4.Range("F22").Value = 1475
